# Fruta / hortaliza, semanal
# The weekly refresh re-pulled the source data, which re-sorted/renumbered the
# daily records for "Terminal Hortofrutícola Agro Chillán - Chirimoya".
# Each data row (2..17) ends up with the Fecha/Calidad/Volumen/Precio* values
# that used to belong to a different row, while the descriptive columns
# (Mercado, Región, Producto, Categoría, Variedad, Unidad, Origen, Kg/unidad)
# stay the same for every row. Capture a snapshot of the old values first,
# then write them back out in their new row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 17

# Snapshot the columns that move between rows: D, L, M, N, O, P, S
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        S = $ws.Cells.Item($r, 19).Value2
    }
}

# Mapping: new row -> old row that supplies its D/L/M/N/O/P/S values
$rowMap = @{
    2  = 12
    3  = 13
    4  = 3
    5  = 16
    6  = 11
    7  = 6
    8  = 4
    9  = 5
    10 = 17
    11 = 14
    12 = 7
    13 = 8
    14 = 15
    15 = 9
    16 = 10
    17 = 2
}

foreach ($newRow in $rowMap.Keys) {
    $oldRow = $rowMap[$newRow]
    $data = $snapshot[$oldRow]

    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 12).Value = $data.L
    $ws.Cells.Item($newRow, 13).Value = $data.M
    $ws.Cells.Item($newRow, 14).Value = $data.N
    $ws.Cells.Item($newRow, 15).Value = $data.O
    $ws.Cells.Item($newRow, 16).Value = $data.P
    $ws.Cells.Item($newRow, 19).Value = $data.S
}
